$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> column -> new value (only cells that actually changed per the diff)
$changes = @{
    2  = @{ D = "67.002.00"; E = "  -3.62%  " }
    3  = @{ D = "3.535.97";  E = "  -3.83%  " }
    4  = @{ D = "1.00";      E = "  -0.07%  " }
    5  = @{ D = "607.01";    E = "  -5.28%  " }
    6  = @{ D = "154.39";    E = "  -3.25%  " }
    7  = @{ D = "3.534.00";  E = "  -3.85%  " }
    8  = @{ E = "  +0.09%  " }
    9  = @{ D = "0.486";     E = "  -2.34%  " }
    10 = @{ E = "  -2.42%  " }
    11 = @{ D = "6.84";      E = "  -3.78%  " }
    12 = @{ E = "  -3.58%  " }
    13 = @{ E = "  -4.39%  " }
    14 = @{ B = "WrappedliquidstakedEther2.0"; C = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"; D = "4.134.26"; E = "  -3.86%  " }
    15 = @{ B = "Avalanche"; C = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"; D = "32.04"; E = "  -2.02%  " }
    16 = @{ D = "3.541.80";  E = "  -3.95%  " }
    17 = @{ D = "66.983.28"; E = "  -3.62%  " }
    18 = @{ E = "  +0.81%  " }
    19 = @{ D = "6.38";      E = "  -1.74%  " }
    20 = @{ D = "15.48" }
    21 = @{ D = "452.14";    E = "  -3.14%  " }
    22 = @{ D = "9.37";      E = "  -5.03%  " }
    23 = @{ E = "  -1.50%  " }
    24 = @{ D = "79.05";     E = "  -0.31%  " }
    25 = @{ D = "3.675.20";  E = "  -3.92%  " }
    26 = @{ E = "  +0.09%  " }
    27 = @{ E = "  -2.00%  " }
    28 = @{ D = "10.25";     E = "  -6.02%  " }
    29 = @{ E = "  -7.51%  " }
    30 = @{ E = "  -0.82%  " }
    31 = @{ E = "  -2.88%  " }
    32 = @{ D = "1.00";      E = "  +0.00%  " }
    33 = @{ E = "  -3.45%  " }
    34 = @{ E = "  -5.15%  " }
    35 = @{ B = "Kaspa"; C = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"; D = "0.158"; E = "  -3.88%  " }
    36 = @{ B = "NEARProtocol"; C = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"; D = "6.19"; E = "  -4.20%  " }
    37 = @{ D = "3.531.75";  E = "  -3.78%  " }
    38 = @{ D = "8.12";      E = "  -4.01%  " }
    39 = @{ E = "  +0.06%  " }
    40 = @{ E = "  -0.09%  " }
    41 = @{ D = "175.90";    E = "  -0.78%  " }
    42 = @{ D = "5.62";      E = "  -4.69%  " }
    43 = @{ E = "  -3.31%  " }
    44 = @{ D = "0.0875";    E = "  -2.67%  " }
    45 = @{ E = "  -3.59%  " }
    46 = @{ D = "45.85";     E = "  -2.18%  " }
    47 = @{ D = "28.39";     E = "  +3.22%  " }
    48 = @{ D = "2.68";      E = "  -1.59%  " }
    49 = @{ E = "  -1.10%  " }
    50 = @{ B = "Cosmos"; C = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"; D = "7.66"; E = "  -2.37%  " }
    51 = @{ B = "SuiNetwork"; C = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"; D = "1.04"; E = "  -3.28%  " }
}

foreach ($row in $changes.Keys) {
    $cols = $changes[$row]
    foreach ($col in $cols.Keys) {
        $newValue = $cols[$col]
        $cell = $ws.Range("$col$row")

        # Values such as "1.00", "6.84" etc. look like numbers and Excel would
        # silently coerce them to the numeric type, losing the original
        # formatted text (e.g. trailing zeros). Force the cell to Text format
        # first so the literal string is preserved, matching the source data
        # which stores every value as text (inline string).
        if ($newValue.Trim() -match '^[0-9]+(\.[0-9]+)?$') {
            $cell.NumberFormat = "@"
        }

        $cell.Value = $newValue
    }
}
